# Add protocol for mining-stone ad NPC configuration:
# 26 new "skeleton_*" NPC rows appended below the existing NPC table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "skeleton_archer_blue",
    "skeleton_archer_green",
    "skeleton_archer_purple",
    "skeleton_archer_red",
    "skeleton_archer_teal",
    "skeleton_archer_yellow",
    "skeleton_tom_angry",
    "skeleton_tom_happy",
    "skeleton_king_blue",
    "skeleton_king_green",
    "skeleton_king_purple",
    "skeleton_king_red",
    "skeleton_king_teal",
    "skeleton_king_yellow",
    "skeleton_mage_blue",
    "skeleton_mage_green",
    "skeleton_mage_purple",
    "skeleton_mage_red",
    "skeleton_mage_teal",
    "skeleton_mage_yellow",
    "skeleton_warrior_blue",
    "skeleton_warrior_green",
    "skeleton_warrior_purple",
    "skeleton_warrior_red",
    "skeleton_warrior_teal",
    "skeleton_warrior_yellow"
)

$startRow = 54
$lastExistingRow = 53

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $startRow + $i
    $name = $names[$i]

    $ws.Range("A$row").Value = $name
    $ws.Range("C$row").Value = "Prefabs/Object/NPC/$name"
    $ws.Range("D$row").Value = 2
    $ws.Range("E$row").Value = 20
    $ws.Range("F$row").Value = "DropBag_1"
    $ws.Range("H$row").Value = 2
    $ws.Range("I$row").Value = "PlayerAtt52"
    $ws.Range("I$row").NumberFormat = $ws.Range("I$lastExistingRow").NumberFormat
    $ws.Range("J$row").Value = "ConsumeData_1"
}

# first new row keeps the running "index" value from the prior row (52)
$ws.Range("B$startRow").Value = 52

$endRow = $startRow + $names.Count - 1

# restore selection/view the author left the sheet in
$ws.Range("J" + $lastExistingRow + ":J" + $endRow).Select()
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 5
